# refactor ADC reading logic
#
# Adds an explicit current (Amps) calculation and a recomputed Rsensor
# sanity-check below the existing Vout/I rows, and inserts one more NTC
# reference data point (100 degC / 6710 Ohm) into the sensor table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 / 16: existing "Vout" / "I (mA)" rows -------------------------
# Promote them visually to the "Calculation" cell style, keeping each
# cell's own number format.
$ws.Range("C15").Style = "Calculation"
$ws.Range("C15").NumberFormat = "0.000"

$ws.Range("C16").Style = "Calculation"
$ws.Range("C16").NumberFormat = "0.00"

# --- Row 17 (new): current through the sensor, in Amps --------------------
$ws.Range("B17").Value = "I"
$ws.Range("C17").Formula = "=(Vin-Vout)/Rref"
$ws.Range("C17").Style = "Calculation"
$ws.Range("C17").NumberFormat = "0.00000"
$ws.Range("D17").Value = "A"

# --- Row 18 (new): Rsensor back-calculated from the measured current ------
$ws.Range("B18").Value = "Rsensor"
$ws.Range("C18").Formula = "=Vout/C17"
$ws.Range("C18").Style = "Calculation"
$ws.Range("C18").NumberFormat = "0"
$ws.Range("D18").Value = "Ohm"

# --- Sensor table: insert a new reference point (100 degC / 6710 Ohm) -----
# before the existing 125 degC row, pushing the rest of the table down one
# row (old row 33 -> 34, old row 34 -> 35).
$ws.Rows.Item(33).Insert()

$ws.Range("C33").Value = 100
$ws.Range("D33").Value = 6710
$ws.Range("E33").Formula = "=Vin/(Rref+D33)*1000"
$ws.Range("F33").Formula = "=E33*1000"
$ws.Range("G33").Formula = "=E33*D33/1000"

# Cursor/selection position, matching the saved workbook state.
$ws.Range("C16").Select()
